# Update Neg_Change and Pos_Change sheets with refreshed filtered market data
$wb = $excel.ActiveWorkbook

# --- Sheet 1: Neg_Change (rows 2-9) ---
$wsNeg = $wb.Worksheets.Item("Neg_Change")

$wsNeg.Cells.Item(2, 1).Value = "IOC"
$wsNeg.Cells.Item(2, 2).Value = 177.28
$wsNeg.Cells.Item(2, 3).Value = 177.9
$wsNeg.Cells.Item(2, 4).Value = 176.04
$wsNeg.Cells.Item(2, 5).Value = 176.72
$wsNeg.Cells.Item(2, 6).Value = 9835632
$wsNeg.Cells.Item(2, 7).Value = 24014261
$wsNeg.Cells.Item(2, 8).Value = -0.5904253726566893
$wsNeg.Cells.Item(2, 9).Value = "IOC"

$wsNeg.Cells.Item(3, 1).Value = "BPCL"
$wsNeg.Cells.Item(3, 2).Value = 377.45
$wsNeg.Cells.Item(3, 3).Value = 378.4
$wsNeg.Cells.Item(3, 4).Value = 372.95
$wsNeg.Cells.Item(3, 5).Value = 374.3
$wsNeg.Cells.Item(3, 6).Value = 2696607
$wsNeg.Cells.Item(3, 7).Value = 6352805
$wsNeg.Cells.Item(3, 8).Value = -0.5755249846327725
$wsNeg.Cells.Item(3, 9).Value = "BPCL"

$wsNeg.Cells.Item(4, 1).Value = "ASTRAL"
$wsNeg.Cells.Item(4, 2).Value = 1589
$wsNeg.Cells.Item(4, 3).Value = 1604.6
$wsNeg.Cells.Item(4, 4).Value = 1559.5
$wsNeg.Cells.Item(4, 5).Value = 1601
$wsNeg.Cells.Item(4, 6).Value = 585825
$wsNeg.Cells.Item(4, 7).Value = 1241019
$wsNeg.Cells.Item(4, 8).Value = -0.527948403690838
$wsNeg.Cells.Item(4, 9).Value = "ASTRAL"

$wsNeg.Cells.Item(5, 1).Value = "HINDPETRO"
$wsNeg.Cells.Item(5, 2).Value = 452.35
$wsNeg.Cells.Item(5, 3).Value = 455.25
$wsNeg.Cells.Item(5, 4).Value = 448.75
$wsNeg.Cells.Item(5, 5).Value = 451
$wsNeg.Cells.Item(5, 6).Value = 1912866
$wsNeg.Cells.Item(5, 7).Value = 3949831
$wsNeg.Cells.Item(5, 8).Value = -0.515709406301181
$wsNeg.Cells.Item(5, 9).Value = "HINDPETRO"

$wsNeg.Cells.Item(6, 1).Value = "SRF"
$wsNeg.Cells.Item(6, 2).Value = 2815
$wsNeg.Cells.Item(6, 3).Value = 2841
$wsNeg.Cells.Item(6, 4).Value = 2765.6
$wsNeg.Cells.Item(6, 5).Value = 2833
$wsNeg.Cells.Item(6, 6).Value = 196550
$wsNeg.Cells.Item(6, 7).Value = 395601
$wsNeg.Cells.Item(6, 8).Value = -0.5031610132431414
$wsNeg.Cells.Item(6, 9).Value = "SRF"

$wsNeg.Cells.Item(7, 1).Value = "FORTIS"
$wsNeg.Cells.Item(7, 2).Value = 928.15
$wsNeg.Cells.Item(7, 3).Value = 930.9
$wsNeg.Cells.Item(7, 4).Value = 914.3
$wsNeg.Cells.Item(7, 5).Value = 919
$wsNeg.Cells.Item(7, 6).Value = 1221319
$wsNeg.Cells.Item(7, 7).Value = 2457834
$wsNeg.Cells.Item(7, 8).Value = -0.5030913397731499
$wsNeg.Cells.Item(7, 9).Value = "FORTIS"

$wsNeg.Cells.Item(8, 1).Value = "SUPREMEIND"
$wsNeg.Cells.Item(8, 2).Value = 3827.9
$wsNeg.Cells.Item(8, 3).Value = 3831.7
$wsNeg.Cells.Item(8, 4).Value = 3745.3
$wsNeg.Cells.Item(8, 5).Value = 3760
$wsNeg.Cells.Item(8, 6).Value = 82044
$wsNeg.Cells.Item(8, 7).Value = 193500
$wsNeg.Cells.Item(8, 8).Value = -0.576
$wsNeg.Cells.Item(8, 9).Value = "SUPREMEIND"

$wsNeg.Cells.Item(9, 1).Value = "HUDCO"
$wsNeg.Cells.Item(9, 2).Value = 201.5
$wsNeg.Cells.Item(9, 3).Value = 201.5
$wsNeg.Cells.Item(9, 4).Value = 193.78
$wsNeg.Cells.Item(9, 5).Value = 195.4
$wsNeg.Cells.Item(9, 6).Value = 10052817
$wsNeg.Cells.Item(9, 7).Value = 22244997
$wsNeg.Cells.Item(9, 8).Value = -0.5480863854465793
$wsNeg.Cells.Item(9, 9).Value = "HUDCO"

# --- Sheet 2: Pos_Change (rows 2-17) ---
$wsPos = $wb.Worksheets.Item("Pos_Change")

$wsPos.Cells.Item(2, 1).Value = "JIOFIN"
$wsPos.Cells.Item(2, 2).Value = 265
$wsPos.Cells.Item(2, 3).Value = 267.25
$wsPos.Cells.Item(2, 4).Value = 262.85
$wsPos.Cells.Item(2, 5).Value = 263.85
$wsPos.Cells.Item(2, 6).Value = 12222314
$wsPos.Cells.Item(2, 7).Value = 8681589
$wsPos.Cells.Item(2, 8).Value = 0.4078429651530383
$wsPos.Cells.Item(2, 9).Value = "JIOFIN"

$wsPos.Cells.Item(3, 1).Value = "HDFCBANK"
$wsPos.Cells.Item(3, 2).Value = 913
$wsPos.Cells.Item(3, 3).Value = 928.2
$wsPos.Cells.Item(3, 4).Value = 901
$wsPos.Cells.Item(3, 5).Value = 903.9
$wsPos.Cells.Item(3, 6).Value = 50188359
$wsPos.Cells.Item(3, 7).Value = 33872203
$wsPos.Cells.Item(3, 8).Value = 0.4816975146257833
$wsPos.Cells.Item(3, 9).Value = "HDFCBANK"

$wsPos.Cells.Item(4, 1).Value = "RECLTD"
$wsPos.Cells.Item(4, 2).Value = 348.1
$wsPos.Cells.Item(4, 3).Value = 348.5
$wsPos.Cells.Item(4, 4).Value = 339.8
$wsPos.Cells.Item(4, 5).Value = 346
$wsPos.Cells.Item(4, 6).Value = 8552636
$wsPos.Cells.Item(4, 7).Value = 5630751
$wsPos.Cells.Item(4, 8).Value = 0.5189156828280987
$wsPos.Cells.Item(4, 9).Value = "RECLTD"

$wsPos.Cells.Item(5, 1).Value = "JSWENERGY"
$wsPos.Cells.Item(5, 2).Value = 479
$wsPos.Cells.Item(5, 3).Value = 479
$wsPos.Cells.Item(5, 4).Value = 468.85
$wsPos.Cells.Item(5, 5).Value = 471.4
$wsPos.Cells.Item(5, 6).Value = 1171877
$wsPos.Cells.Item(5, 7).Value = 831728
$wsPos.Cells.Item(5, 8).Value = 0.4089666333224323
$wsPos.Cells.Item(5, 9).Value = "JSWENERGY"

$wsPos.Cells.Item(6, 1).Value = "BOSCHLTD"
$wsPos.Cells.Item(6, 2).Value = 36415
$wsPos.Cells.Item(6, 3).Value = 36420
$wsPos.Cells.Item(6, 4).Value = 35600
$wsPos.Cells.Item(6, 5).Value = 35700
$wsPos.Cells.Item(6, 6).Value = 13848
$wsPos.Cells.Item(6, 7).Value = 9373
$wsPos.Cells.Item(6, 8).Value = 0.4774351861730503
$wsPos.Cells.Item(6, 9).Value = "BOSCHLTD"

$wsPos.Cells.Item(7, 1).Value = "SOLARINDS"
$wsPos.Cells.Item(7, 2).Value = 13250
$wsPos.Cells.Item(7, 3).Value = 13300
$wsPos.Cells.Item(7, 4).Value = 12935
$wsPos.Cells.Item(7, 5).Value = 12965
$wsPos.Cells.Item(7, 6).Value = 90669
$wsPos.Cells.Item(7, 7).Value = 60892
$wsPos.Cells.Item(7, 8).Value = 0.4890133350850687
$wsPos.Cells.Item(7, 9).Value = "SOLARINDS"

$wsPos.Cells.Item(8, 1).Value = "DLF"
$wsPos.Cells.Item(8, 2).Value = 641.6
$wsPos.Cells.Item(8, 3).Value = 641.6
$wsPos.Cells.Item(8, 4).Value = 623.55
$wsPos.Cells.Item(8, 5).Value = 626.95
$wsPos.Cells.Item(8, 6).Value = 4838136
$wsPos.Cells.Item(8, 7).Value = 3109096
$wsPos.Cells.Item(8, 8).Value = 0.5561230659973189
$wsPos.Cells.Item(8, 9).Value = "DLF"

$wsPos.Cells.Item(9, 1).Value = "AUROPHARMA"
$wsPos.Cells.Item(9, 2).Value = 1149.1
$wsPos.Cells.Item(9, 3).Value = 1157.4
$wsPos.Cells.Item(9, 4).Value = 1126.8
$wsPos.Cells.Item(9, 5).Value = 1157
$wsPos.Cells.Item(9, 6).Value = 1682470
$wsPos.Cells.Item(9, 7).Value = 1090045
$wsPos.Cells.Item(9, 8).Value = 0.5434867367860959
$wsPos.Cells.Item(9, 9).Value = "AUROPHARMA"

$wsPos.Cells.Item(10, 1).Value = "CUMMINSIND"
$wsPos.Cells.Item(10, 2).Value = 4425
$wsPos.Cells.Item(10, 3).Value = 4491.9
$wsPos.Cells.Item(10, 4).Value = 4395.2
$wsPos.Cells.Item(10, 5).Value = 4412
$wsPos.Cells.Item(10, 6).Value = 512643
$wsPos.Cells.Item(10, 7).Value = 360796
$wsPos.Cells.Item(10, 8).Value = 0.4208666393197264
$wsPos.Cells.Item(10, 9).Value = "CUMMINSIND"

$wsPos.Cells.Item(11, 1).Value = "MPHASIS"
$wsPos.Cells.Item(11, 2).Value = 2387.6
$wsPos.Cells.Item(11, 3).Value = 2501.2
$wsPos.Cells.Item(11, 4).Value = 2355.2
$wsPos.Cells.Item(11, 5).Value = 2452.5
$wsPos.Cells.Item(11, 6).Value = 1725163
$wsPos.Cells.Item(11, 7).Value = 1136104
$wsPos.Cells.Item(11, 8).Value = 0.5184903846830924
$wsPos.Cells.Item(11, 9).Value = "MPHASIS"

$wsPos.Cells.Item(12, 1).Value = "LUPIN"
$wsPos.Cells.Item(12, 2).Value = 2210
$wsPos.Cells.Item(12, 3).Value = 2219.2
$wsPos.Cells.Item(12, 4).Value = 2160
$wsPos.Cells.Item(12, 5).Value = 2194
$wsPos.Cells.Item(12, 6).Value = 1126141
$wsPos.Cells.Item(12, 7).Value = 804229
$wsPos.Cells.Item(12, 8).Value = 0.4002740512963348
$wsPos.Cells.Item(12, 9).Value = "LUPIN"

$wsPos.Cells.Item(13, 1).Value = "DIXON"
$wsPos.Cells.Item(13, 2).Value = 11500
$wsPos.Cells.Item(13, 3).Value = 11548
$wsPos.Cells.Item(13, 4).Value = 11250
$wsPos.Cells.Item(13, 5).Value = 11398
$wsPos.Cells.Item(13, 6).Value = 446533
$wsPos.Cells.Item(13, 7).Value = 288312
$wsPos.Cells.Item(13, 8).Value = 0.5487839562695969
$wsPos.Cells.Item(13, 9).Value = "DIXON"

$wsPos.Cells.Item(14, 1).Value = "BANKINDIA"
$wsPos.Cells.Item(14, 2).Value = 164.1
$wsPos.Cells.Item(14, 3).Value = 164.45
$wsPos.Cells.Item(14, 4).Value = 161.26
$wsPos.Cells.Item(14, 5).Value = 161.69
$wsPos.Cells.Item(14, 6).Value = 8120034
$wsPos.Cells.Item(14, 7).Value = 5243794
$wsPos.Cells.Item(14, 8).Value = 0.5485036216144265
$wsPos.Cells.Item(14, 9).Value = "BANKINDIA"

$wsPos.Cells.Item(15, 1).Value = "BSE"
$wsPos.Cells.Item(15, 2).Value = 3100
$wsPos.Cells.Item(15, 3).Value = 3127
$wsPos.Cells.Item(15, 4).Value = 3017.9
$wsPos.Cells.Item(15, 5).Value = 3023.9
$wsPos.Cells.Item(15, 6).Value = 3885000
$wsPos.Cells.Item(15, 7).Value = 2525871
$wsPos.Cells.Item(15, 8).Value = 0.5380832987907933
$wsPos.Cells.Item(15, 9).Value = "BSE"

$wsPos.Cells.Item(16, 1).Value = "CROMPTON"
$wsPos.Cells.Item(16, 2).Value = 264
$wsPos.Cells.Item(16, 3).Value = 269.85
$wsPos.Cells.Item(16, 4).Value = 255.9
$wsPos.Cells.Item(16, 5).Value = 266.8
$wsPos.Cells.Item(16, 6).Value = 4504154
$wsPos.Cells.Item(16, 7).Value = 3035533
$wsPos.Cells.Item(16, 8).Value = 0.4838099272845988
$wsPos.Cells.Item(16, 9).Value = "CROMPTON"

$wsPos.Cells.Item(17, 1).Value = "CDSL"
$wsPos.Cells.Item(17, 2).Value = 1350
$wsPos.Cells.Item(17, 3).Value = 1359
$wsPos.Cells.Item(17, 4).Value = 1321.1
$wsPos.Cells.Item(17, 5).Value = 1335.6
$wsPos.Cells.Item(17, 6).Value = 1663043
$wsPos.Cells.Item(17, 7).Value = 1133736
$wsPos.Cells.Item(17, 8).Value = 0.4668697121728515
$wsPos.Cells.Item(17, 9).Value = "CDSL"
